$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D and G to remain text so numeric-looking strings are preserved exactly
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# Update Price (column D) values that changed
$ws.Range("D2").Value = "244.52"
$ws.Range("D3").Value = "24.08"
$ws.Range("D4").Value = "5.858"
$ws.Range("D5").Value = "0.05856"
$ws.Range("D6").Value = "3.428"
$ws.Range("D7").Value = "6.505"
$ws.Range("D8").Value = "1.330"
$ws.Range("D9").Value = "0.7988"
$ws.Range("D10").Value = "0.1473"
$ws.Range("D11").Value = "0.07717"
$ws.Range("D12").Value = "0.03292"
$ws.Range("D13").Value = "0.03015"
$ws.Range("D14").Value = "0.09213"
$ws.Range("D15").Value = "3.561"
$ws.Range("D16").Value = "0.001669"
$ws.Range("D17").Value = "0.04773"
$ws.Range("D18").Value = "0.0006039"
$ws.Range("D19").Value = "0.006237"
$ws.Range("D20").Value = "0.005516"
$ws.Range("D21").Value = "0.001070"
$ws.Range("D22").Value = "0.0001502"
$ws.Range("D23").Value = "3.703"
$ws.Range("D25").Value = "0.3329"
$ws.Range("D26").Value = "0.1253"
$ws.Range("D27").Value = "0.0006279"
$ws.Range("D40").Value = "0.04357"
$ws.Range("D41").Value = "0.007068"
$ws.Range("D42").Value = "0.003604"
$ws.Range("D43").Value = "0.1062"
$ws.Range("D44").Value = "0.008700"
$ws.Range("D45").Value = "0.002464"
$ws.Range("D46").Value = "0.00005897"
$ws.Range("D48").Value = "0.9914"
$ws.Range("D49").Value = "0.1087"
$ws.Range("D50").Value = "0.00002103"

# Update Hora (column G) values: all rows 2-51 changed from 15 to 16
$ws.Range("G2:G51").Value = "16"

